$d = $word.ActiveDocument

# Locate the "#Ubuntu/Debian" run and collapse the range to its end point.
$rng = $d.Content
$rng.Find.Execute("#Ubuntu/Debian", $false, $true, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)

# Insert a new paragraph right after "#Ubuntu/Debian" (before the bookmark).
$rng.InsertParagraphAfter()

# Move into the newly created (still empty) paragraph and set its text.
$rng.MoveStart(4, 1)
$rng.MoveEnd(4, 1)
$rng.Text = "Add a file called 30-phalcon.ini in /etc/php.d/ with this content: extension=phalcon.so"
